{"js": "// Update the BCA -> MCA wording and the academic year (2014-2015 -> 2017-2018)\n// across the certificate document.\n\nconst body = context.document.body;\n\nasync function replaceOnce(searchText, newText) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, \"Replace\");\n    await context.sync();\n  }\n}\n\n// \"B.C.A.III (Semester-VI)\" -> \"M.C.A (Semester-VI)\"\nawait replaceOnce(\"B.C.A.III\", \"M.C.A\");\n\n// \"(2014-2015)\" -> \"(2017-2018)\"\nawait replaceOnce(\"(2014-2015)\", \"(2017-2018)\");\n\n// \"...MUMBAI (SESCOM   BCA)\" -> \"...MUMBAI (SESCOM   MCA)\"\nawait replaceOnce(\" BCA)\", \" MCA)\");\n\n// \"...academic year 2014-2015 by\" -> \"...academic year 2017-2018 by\"\nawait replaceOnce(\"2014\", \"2017\");\nawait replaceOnce(\"2015\", \"2018\");\n", "ps1": "# Update the BCA -> MCA wording and the academic year (2014-2015 -> 2017-2018)\n# across the certificate document.\n\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n# \"B.C.A.III (Semester-VI)\" -> \"M.C.A (Semester-VI)\"\n$d.Content.Find.Execute(\"B.C.A.III\", $false, $true, $false, $false, $false, $true, `\n    $wdFindContinue, $false, \"M.C.A\", $wdReplaceAll) | Out-Null\n\n# \"(2014-2015)\" -> \"(2017-2018)\"\n$d.Content.Find.Execute(\"(2014-2015)\", $false, $true, $false, $false, $false, $true, `\n    $wdFindContinue, $false, \"(2017-2018)\", $wdReplaceAll) | Out-Null\n\n# \"...MUMBAI (SESCOM   BCA)\" -> \"...MUMBAI (SESCOM   MCA)\"\n$d.Content.Find.Execute(\" BCA)\", $false, $true, $false, $false, $false, $true, `\n    $wdFindContinue, $false, \" MCA)\", $wdReplaceAll) | Out-Null\n\n# \"...academic year 2014-2015 by\" -> \"...academic year 2017-2018 by\"\n$d.Content.Find.Execute(\"2014\", $false, $true, $false, $false, $false, $true, `\n    $wdFindContinue, $false, \"2017\", $wdReplaceAll) | Out-Null\n$d.Content.Find.Execute(\"2015\", $false, $true, $false, $false, $false, $true, `\n    $wdFindContinue, $false, \"2018\", $wdReplaceAll) | Out-Null\n"}
